$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Kanban sheet: move three cards to new columns/rows and add a new one.
#   - "API para uso externo com resposta json" moves from "To do" (A17) to
#     "Done" (C16) and its effort tag drops from (****) to (**).
#   - "Documentação do aplicativo(***)" moves from "W.I.P." (B16) to
#     "Done" (C17).
#   - "Correção de Issues e melhorias(****)" moves from "To do" (A18) to
#     "W.I.P." (B18).
#   - "Subir para produção (****) (****)" stays put in "To do" (A19).
# ---------------------------------------------------------------------------
$kanban = $wb.Worksheets.Item("Kanban")

$kanban.Range("B16").ClearContents()
$kanban.Range("C16").Value = "API para uso externo com resposta json (**)"

$kanban.Range("A17").ClearContents()
$kanban.Range("C17").Value = "Documentação do aplicativo(***)"

$kanban.Range("A18").ClearContents()
$kanban.Range("B18").Value = "Correção de Issues e melhorias(****)"

[void]$kanban.Range("A18").Select()

# ---------------------------------------------------------------------------
# Backlog sheet: mark the matching backlog items as "Done" (same fill/style
# already used by rows 2-11, the green "Bom" cell style).
# ---------------------------------------------------------------------------
$backlog = $wb.Worksheets.Item("Backlog")
$backlog.Range("A12:C13").Style = "Bom"

# Restore "Backlog" as the active sheet (selecting on Kanban above made it
# active as a side effect).
[void]$backlog.Activate()
